$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Shift the existing "Language"/english/spanish column (old column C)
#     out to the new column E, before we overwrite column C with new data.
$ws.Range("C1:C3").Copy()
$ws.Range("E1:E3").PasteSpecial(-4104)

# --- Row 1 (header) ---
$ws.Range("B1").Value = "NavGroup"
$ws.Range("C1").Value = "SubNavGroup"
$ws.Range("D1").Value = "ListItem"

# Copy the header style (bold font + gray fill, style index 1 from A1) onto
# the newly-populated header cells so they match the rest of the row.
$ws.Range("A1").Copy()
$ws.Range("B1:E1").PasteSpecial(-4122)

# --- Row 3 ---
$ws.Range("B3").Value = "Tipos de cáncer"
$ws.Range("C3").Value = "Tipos comunes de cáncer"
$ws.Range("D3").Value = "Linfoma"

# --- Row 2 ---
$ws.Range("B2").Value = "About Cancer"
$ws.Range("C2").Value = "Understanding Cancer"
$ws.Range("D2").Value = "Cancer Statistics"

# --- Column widths: C gets a bestFit width, D:E match column B's width ---
$ws.Columns.Item(3).ColumnWidth = 22.7
$ws.Columns.Item(4).ColumnWidth = 23.3
$ws.Columns.Item(5).ColumnWidth = 23.3
